$wb = $excel.ActiveWorkbook

# --- Sheet "全部类型" (sheet4): row 14 restructuring ---
# A new event ("上海·东方明珠...") was inserted as row 14, the old row-14
# ("X-party") content shifts down to row 15 (values unchanged), the old
# row-15 ("BH盛会之星") content shifts down to row 16 (F value bumped), and
# the old row-16 ("VWonderland") entry is removed entirely. Row 17 onward
# keeps its original identity (only the "want to go" counter changes).
$ws4 = $wb.Worksheets.Item("全部类型")

# New row 14 content (replaces the old X-party row in place)
$ws4.Range("B14").Value = '2024-08-17'
$ws4.Range("C14").Value = '上海·东方明珠·「光与夜之恋 × 线条小狗 ×爱胖达文化 」线条大作战主题店'
$ws4.Range("D14").Value = '世纪大道1号 东方明珠电视塔城市广场商场'
$ws4.Range("E14").Value = '2024.08.17 00:00-10.27 23:59'
$ws4.Range("F14").Value = 1853
$ws4.Range("G14").Value = '已售罄'
$ws4.Range("H14").Value = 'https://show.bilibili.com/platform/detail.html?id=90444'
$ws4.Range("I14").Value = '//i2.hdslb.com/bfs/openplatform/202408/qUE9n4UR1723020534077.png'

# Row 15 now holds the content that used to be in row 14 (X-party) - values unchanged
$ws4.Range("C15").Value = ' 上海·X-party国漫游戏嘉年华02（免费活动）'
$ws4.Range("D15").Value = '七莘路1599弄(七莘路地铁站1号口步行60米) 平金中心'
$ws4.Range("E15").Value = '2024.08.24 13:00-08.24 18:00'
$ws4.Range("F15").Value = 38
$ws4.Range("G15").Value = 48
$ws4.Range("H15").Value = 'https://show.bilibili.com/platform/detail.html?id=90885'
$ws4.Range("I15").Value = '//i1.hdslb.com/bfs/openplatform/202408/cNzKejgM1724147367658.jpeg'

# Row 16 now holds the content that used to be in row 15 (BH盛会之星) - F value +1
$ws4.Range("C16").Value = '上海·BH盛会之星的邀约'
$ws4.Range("D16").Value = '鲁班路300号 星光摄影器材城'
$ws4.Range("E16").Value = '2024.08.24 10:30-08.25 17:00'
$ws4.Range("F16").Value = 439
$ws4.Range("G16").Value = 60
$ws4.Range("H16").Value = 'https://show.bilibili.com/platform/detail.html?id=88603'
$ws4.Range("I16").Value = '//i1.hdslb.com/bfs/openplatform/202407/fUi7Oz2b1719995931315.png'
# (The old row-16 "VWonderland" entry is removed by being overwritten above;
#  rows 17 onward keep their original identity.)

# --- Sheet "展览" (sheet1): update "want to go" counts (column F) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3145
$ws.Range("F9").Value = 7373
$ws.Range("F11").Value = 149
$ws.Range("F13").Value = 439
$ws.Range("F15").Value = 1942
$ws.Range("F16").Value = 1797
$ws.Range("F17").Value = 1090
$ws.Range("F20").Value = 1853
$ws.Range("F21").Value = 1386
$ws.Range("F22").Value = 1260
$ws.Range("F23").Value = 657
$ws.Range("F25").Value = 1144
$ws.Range("F27").Value = 547
$ws.Range("F29").Value = 76
$ws.Range("F30").Value = 4717
$ws.Range("F32").Value = 3918
$ws.Range("F33").Value = 2197
$ws.Range("F34").Value = 166
$ws.Range("F35").Value = 228
$ws.Range("F40").Value = 385
$ws.Range("F41").Value = 1
$ws.Range("F42").Value = 160
$ws.Range("F43").Value = 526
$ws.Range("F44").Value = 260
$ws.Range("F45").Value = 201
$ws.Range("F46").Value = 804
$ws.Range("F47").Value = 437
$ws.Range("F48").Value = 10
$ws.Range("F49").Value = 160

# --- Sheet "演出" (sheet2): update "want to go" counts (column F) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 96
$ws.Range("F24").Value = 94

# --- Sheet "本地生活" (sheet3): update "want to go" counts (column F) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 565
$ws.Range("F6").Value = 1855
$ws.Range("F8").Value = 2909
$ws.Range("F9").Value = 1147
$ws.Range("F12").Value = 440
$ws.Range("F13").Value = 1853
$ws.Range("F14").Value = 8225
$ws.Range("F15").Value = 400

# --- Sheet "全部类型" (sheet4): remaining "want to go" count updates (column F) ---
$ws4.Range("F4").Value = 3145
$ws4.Range("F6").Value = 1855
$ws4.Range("F8").Value = 2909
$ws4.Range("F9").Value = 7373
$ws4.Range("F10").Value = 1147
$ws4.Range("F13").Value = 440
$ws4.Range("F17").Value = 1090
$ws4.Range("F20").Value = 1853
$ws4.Range("F21").Value = 1386
$ws4.Range("F22").Value = 1260
$ws4.Range("F23").Value = 657
$ws4.Range("F25").Value = 1144
$ws4.Range("F26").Value = 96
$ws4.Range("F30").Value = 547
$ws4.Range("F33").Value = 76
$ws4.Range("F34").Value = 4718
$ws4.Range("F36").Value = 3918
$ws4.Range("F37").Value = 2197
$ws4.Range("F38").Value = 166
$ws4.Range("F39").Value = 228
$ws4.Range("F43").Value = 385
$ws4.Range("F44").Value = 160
$ws4.Range("F45").Value = 94
$ws4.Range("F46").Value = 526
$ws4.Range("F47").Value = 260
$ws4.Range("F49").Value = 437

Write-Host "done"
